# StoryCards.xlsx — add three new story-card entries (rows 14, 15, 16 / sheet
# rows 21-23) to "Tabelle1", matching the commit
# "Erstellen von neuen Einträgen in StoryCards.xlsx."

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# New data rows (sheet rows 21-23), mirroring the layout of the existing
# rows 8-20: A=ID, B=Status, C=Priorität, D=Art, E=Datum, F=Autor,
# G=Entwickler, H=Titel, I=Beschreibung (J..M left blank, same as rows 14/17/20)
$newRows = @(
    @{ Row = 21; A = 14; B = "jungfräulich"; C = "mittel"; D = "Öffentlichkeitsarbeit"; E = 40823; F = "Salzer"; G = "Wiederschein"; H = "Recherche";    I = "Sammeln von infos zu Ameisensysteme" },
    @{ Row = 22; A = 15; B = "jungfräulich"; C = "mittel"; D = "Öffentlichkeitsarbeit"; E = 40823; F = "Salzer"; G = "Wiederschein"; H = "Recherche";    I = "Sammeln von infos zu TSP" },
    @{ Row = 23; A = 16; B = "jungfräulich"; C = "hoch";   D = "Öffentlichkeitsarbeit"; E = 40823; F = "Salzer"; G = "Wiederschein"; H = "Ausarbeitung"; I = "Komprimierung der gesammelten Daten" }
)

foreach ($entry in $newRows) {
    $r = $entry.Row

    $ws.Cells.Item($r, 1).Value = $entry.A
    $ws.Cells.Item($r, 2).Value = $entry.B
    $ws.Cells.Item($r, 3).Value = $entry.C
    $ws.Cells.Item($r, 4).Value = $entry.D

    # Date column: write the serial value first, then copy the existing
    # date cell's number format (style index) across without touching the
    # value, so it reuses the workbook's existing date style instead of
    # creating a brand-new number format.
    $ws.Cells.Item($r, 5).Value = $entry.E
    $ws.Cells.Item(8, 5).Copy() | Out-Null
    $ws.Cells.Item($r, 5).PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($r, 6).Value = $entry.F
    $ws.Cells.Item($r, 7).Value = $entry.G
    $ws.Cells.Item($r, 8).Value = $entry.H
    $ws.Cells.Item($r, 9).Value = $entry.I
}

$excel.CutCopyMode = 0

# Match the saved selection/active cell state (the sheet view no longer
# scrolls to keep column B at the left edge, and the last selection is G26).
$ws.Activate()
$ws.Range("G26").Select() | Out-Null
